$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44162
$ws.Range("H2").Value = 'Sin especificar'
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17500
$ws.Range("N2").Value = '$/saco 25 kilos'
$ws.Range("O2").Value = 'Región del Maule'
$ws.Range("P2").Value = 700

# Row 3
$ws.Range("D3").Value = 44454
$ws.Range("H3").Value = 'Perfection'
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 36000
$ws.Range("L3").Value = 38000
$ws.Range("M3").Value = 37000
$ws.Range("N3").Value = '$/malla 25 kilos'
$ws.Range("O3").Value = 'Provincia de Limarí'
$ws.Range("P3").Value = 1480

# Row 4
$ws.Range("D4").Value = 44482
$ws.Range("H4").Value = 'Perfection'
$ws.Range("J4").Value = 130
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 24385
$ws.Range("N4").Value = '$/saco 25 kilos'
$ws.Range("O4").Value = 'Región de O''Higgins'
$ws.Range("P4").Value = 975

# Row 5
$ws.Range("D5").Value = 44483
$ws.Range("H5").Value = 'Perfection'
$ws.Range("J5").Value = 220
$ws.Range("K5").Value = 19000
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = 19455
$ws.Range("N5").Value = '$/saco 25 kilos'
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 778

# Row 6
$ws.Range("D6").Value = 44517
$ws.Range("H6").Value = 'Perfection'
$ws.Range("J6").Value = 110
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17455
$ws.Range("N6").Value = '$/saco 25 kilos'
$ws.Range("O6").Value = 'Región del Maule'
$ws.Range("P6").Value = 698

# Row 7
$ws.Range("D7").Value = 44328
$ws.Range("H7").Value = 'Perfection'
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 33000
$ws.Range("L7").Value = 34000
$ws.Range("M7").Value = 33500
$ws.Range("N7").Value = '$/malla 25 kilos'
$ws.Range("O7").Value = 'Provincia de Huasco'
$ws.Range("P7").Value = 1340

# Row 8
$ws.Range("D8").Value = 44503
$ws.Range("H8").Value = 'Perfection'
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("N8").Value = '$/malla 25 kilos'
$ws.Range("O8").Value = 'Provincia de Limarí'
$ws.Range("P8").Value = 620

# Row 9
$ws.Range("D9").Value = 44505
$ws.Range("H9").Value = 'Perfection'
$ws.Range("J9").Value = 210
$ws.Range("K9").Value = 6500
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 6714
$ws.Range("N9").Value = '$/malla 25 kilos'
$ws.Range("O9").Value = 'Región del Maule'
$ws.Range("P9").Value = 269

# Row 10
$ws.Range("D10").Value = 44496
$ws.Range("H10").Value = 'Perfection'
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14520
$ws.Range("N10").Value = '$/malla 25 kilos'
$ws.Range("O10").Value = 'Provincia de Huasco'
$ws.Range("P10").Value = 581

# Row 11
$ws.Range("D11").Value = 44518
$ws.Range("H11").Value = 'Perfection'
$ws.Range("J11").Value = 350
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14571
$ws.Range("N11").Value = '$/saco 25 kilos'
$ws.Range("O11").Value = 'Región del Maule'
$ws.Range("P11").Value = 583

# Row 12
$ws.Range("D12").Value = 44519
$ws.Range("H12").Value = 'Perfection'
$ws.Range("J12").Value = 240
$ws.Range("K12").Value = 17000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 17583
$ws.Range("N12").Value = '$/saco 25 kilos'
$ws.Range("O12").Value = 'Carahue'
$ws.Range("P12").Value = 703

# Row 13
$ws.Range("D13").Value = 44399
$ws.Range("H13").Value = 'Perfection'
$ws.Range("J13").Value = 50
$ws.Range("K13").Value = 39000
$ws.Range("L13").Value = 40000
$ws.Range("M13").Value = 39600
$ws.Range("N13").Value = '$/malla 25 kilos'
$ws.Range("O13").Value = 'Provincia de Huasco'
$ws.Range("P13").Value = 1584

# Row 14
$ws.Range("D14").Value = 44512
$ws.Range("H14").Value = 'Perfection'
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14500
$ws.Range("N14").Value = '$/saco 25 kilos'
$ws.Range("O14").Value = 'Región del Maule'
$ws.Range("P14").Value = 580

# Row 15
$ws.Range("D15").Value = 44335
$ws.Range("H15").Value = 'Perfection'
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 30000
$ws.Range("L15").Value = 32000
$ws.Range("M15").Value = 31000
$ws.Range("N15").Value = '$/malla 25 kilos'
$ws.Range("O15").Value = 'Provincia de Huasco'
$ws.Range("P15").Value = 1240

# Row 16
$ws.Range("D16").Value = 44342
$ws.Range("H16").Value = 'Perfection'
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 30000
$ws.Range("L16").Value = 32000
$ws.Range("M16").Value = 31000
$ws.Range("N16").Value = '$/malla 25 kilos'
$ws.Range("O16").Value = 'Provincia de Limarí'
$ws.Range("P16").Value = 1240

